$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in column C
$ws.Range("C3").Value = 51
$ws.Range("C12").Value = -56.4
$ws.Range("C53").Value = -57.4

# Update the active cell selection on the sheet
$ws.Range("F8").Select()
